# Update the "adults" worksheet with the revised COVID-19 dosing rules and
# the new "Atleast" wording (typo fix from "At lest") plus "spacing" suffixes
# on the interval columns. Also add a new COVID-19 bivalent-booster row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("adults")

# --- Row 2 (COVID-19 primary series -> now 1 dose bivalent) ---
$ws.Range("B2").Value = 1
$ws.Range("C2").Value = "Atleast 19 years old -> 1 dose bivalent"
$ws.Range("D2").Value = "X"
$ws.Range("H2").Value = 6935
$ws.Range("I2").Value = 23360
$ws.Range("K2").Value = 23360

# --- Row 3 (Influenza) ---
$ws.Range("C3").Value = "Atleast 19 years old - Annually "

# --- Row 4 (Tdap) ---
$ws.Range("C4").Value = "Atleast 11 years old"

# --- Row 5 (MMR) ---
$ws.Range("C5").Value = "Atleast 19 years old"

# --- Row 6 (Varicella) ---
$ws.Range("C6").Value = "Atleast 19 years old"
$ws.Range("D6").Value = "4 weeks spacing"

# --- Row 7 (Zoster / Shingles) ---
$ws.Range("C7").Value = "Atleast 50 years old"
$ws.Range("D7").Value = "2 months spacing"

# --- Row 8 (HPV) ---
$ws.Range("C8").Value = "Atleast 9 years old"
$ws.Range("D8").Value = "5 months spacing if first dose was between 9-14 years old, 1 month spacing if first dose was 15 years or older"
$ws.Range("E8").Value = "6 months spacing (only if older than 15 years old for first dose)"

# --- Row 9 (Pneumococcal) ---
$ws.Range("C9").Value = "Atleast 65 years old"
$ws.Range("D9").Value = "1 year spacing (if PCV15 used)"

# --- Row 10 (Hepatitis B) ---
$ws.Range("C10").Value = "Atleast 19 years old"
$ws.Range("D10").Value = "1 month spacing"
$ws.Range("E10").Value = "6 months spacing"

# --- New Row 11: COVID-19 bivalent booster entry ---
# Clone row 2's formatting (cell styles) and values, then adjust the
# vaccine-specific fields for the new entry.
$ws.Range("A2:AW2").Copy()
$ws.Range("A11").PasteSpecial(-4122)
$ws.Range("A2:AW2").Copy()
$ws.Range("A11").PasteSpecial(-4163)
$excel.CutCopyMode = 0
$ws.Range("N11:S11").Clear()

$ws.Range("A11").Value = "COVID-19"
$ws.Range("B11").Value = 2
$ws.Range("C11").Value = "Atleast 65 years old -> 2 doses bivalent"
$ws.Range("D11").Value = "4 months spacing"
$ws.Range("E11").Value = "X"
$ws.Range("F11").Value = "X"
$ws.Range("G11").Value = "X"
$ws.Range("H11").Value = 23361
$ws.Range("I11").Value = 364635
$ws.Range("J11").Value = 23361
$ws.Range("K11").Value = 364635
$ws.Range("L11").Value = 23361
$ws.Range("M11").Value = 364635

# --- Move the selection (view now scrolls back to show column A) ---
$ws.Range("L23").Select()
